$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Insert 5 new columns (E:I) shifting the existing E:I data to J:N
$ws.Range("E1:I31").Insert(-4161)

# --- Header row 8 (first table) ---
$ws.Range("E8").Value = "فصل دوم منتهی به 1398/12"
$ws.Range("F8").Value = "فصل سوم منتهی به 1399/03"
$ws.Range("G8").Value = "فصل چهارم منتهی به 1399/06"
$ws.Range("H8").Value = "فصل اول منتهی به 1399/09"
$ws.Range("I8").Value = "فصل دوم منتهی به 1399/12"

# --- Header row 24 (second table) ---
$ws.Range("E24").Value = "فصل دوم منتهی به 1398/12"
$ws.Range("F24").Value = "فصل سوم منتهی به 1399/03"
$ws.Range("G24").Value = "فصل چهارم منتهی به 1399/06"
$ws.Range("H24").Value = "فصل اول منتهی به 1399/09"
$ws.Range("I24").Value = "فصل دوم منتهی به 1399/12"

# --- Row 10: هزینه حمل و نقل و انتقال ---
$ws.Range("E10").Value = -25106
$ws.Range("F10").Value = 176085
$ws.Range("G10").Value = -87368
$ws.Range("H10").Value = 49208
$ws.Range("I10").Value = 29466

# --- Row 11: هزینه خدمات پس از فروش ---
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0

# --- Row 12: حق العمل و کمیسیون فروش ---
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0

# --- Row 13: هزینه تبلیغات ---
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0

# --- Row 14: هزینه مواد مصرفی ---
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0

# --- Row 15: هزینه انرژی (آب، برق، گاز و سوخت) ---
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0

# --- Row 16: هزینه استهلاک ---
$ws.Range("E16").Value = 6710
$ws.Range("F16").Value = 6985
$ws.Range("G16").Value = 6925
$ws.Range("H16").Value = 6656
$ws.Range("I16").Value = 5977

# --- Row 17: هزینه حقوق و دستمزد ---
$ws.Range("E17").Value = 31988
$ws.Range("F17").Value = 34859
$ws.Range("G17").Value = 56087
$ws.Range("H17").Value = 48308
$ws.Range("I17").Value = 57804

# --- Row 18: هزینه مطالبات مشکوک الوصول ---
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0

# --- Row 19: سایر هزینه ها ---
$ws.Range("E19").Value = 706278
$ws.Range("F19").Value = 383153
$ws.Range("G19").Value = 6177223
$ws.Range("H19").Value = 3727051
$ws.Range("I19").Value = 3732080

# --- Row 20: جمع ---
$ws.Range("E20").Value = 719870
$ws.Range("F20").Value = 601082
$ws.Range("G20").Value = 6152867
$ws.Range("H20").Value = 3831223
$ws.Range("I20").Value = 3825327

# --- Row 26: تعداد پرسنل غیر تولیدی شرکت ---
$ws.Range("E26").Value = 180
$ws.Range("F26").Value = 181
$ws.Range("G26").Value = 181
$ws.Range("H26").Value = 143
$ws.Range("I26").Value = 172

# --- Row 27: تعداد پرسنل تولیدی شرکت ---
$ws.Range("E27").Value = 513
$ws.Range("F27").Value = 510
$ws.Range("G27").Value = 511
$ws.Range("H27").Value = 549
$ws.Range("I27").Value = 520
